$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old local path prefix with the new network-drive path prefix
# across all data rows (2-61) in column C.
$oldPrefix = "C:/Users/Pedro/Documents/Covid19VaccinationProject"
$newPrefix = "W:/Meu Drive/Doutorado IC-Unicamp/MO412 - Graphs Algorithms/07. Final Project/Covid19VaccinationProject"

for ($row = 2; $row -le 61; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value()
    $cell.Value = $current.Replace($oldPrefix, $newPrefix)
}

# Update the recorded row counts (column D) that changed as a result of the
# re-processing against the new (larger/updated) source files.
$rowCounts = @{
    2 = 354957
    3 = 1077268
    4 = 236189
    5 = 1026741
    6 = 5109864
    7 = 2660820
    8 = 103065
    9 = 1252587
    10 = 1973666
    11 = 2557270
    12 = 1301109
    13 = 979687
    14 = 7097066
    15 = 2323741
    16 = 1585074
    17 = 3361230
    18 = 3085171
    19 = 1314846
    20 = 2124822
    21 = 1530749
    22 = 3724089
    23 = 539053
    24 = 199068
    25 = 2520609
    26 = 9255253
    27 = 868964
    28 = 748616
    29 = 354957
    30 = 1077268
    31 = 236189
    32 = 1026741
    33 = 5109864
    34 = 2660820
    35 = 103065
    36 = 1252587
    37 = 1973666
    38 = 2557270
    39 = 1301109
    40 = 979687
    41 = 7097066
    42 = 2323741
    43 = 1585074
    44 = 3361230
    45 = 3085171
    46 = 1314846
    47 = 2124822
    48 = 1530749
    49 = 3724089
    50 = 539053
    51 = 199068
    52 = 2520609
    53 = 9255253
    54 = 868964
    55 = 748616
    58 = 1751861
    61 = 1751861
}

foreach ($row in $rowCounts.Keys) {
    $ws.Cells.Item($row, 4).Value = $rowCounts[$row]
}

